$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Assigning a plain numeric-looking string via .Value lets Excel
    # auto-convert it to a number, which would lose the original text
    # formatting (leading zeros, trailing zeros, thousand-dot grouping).
    # Force it in as a formula that evaluates to the literal text, then
    # copy/paste-special as values to bake it into a plain text cell
    # without leaving the formula behind or touching the cell style.
    $rng = $ws.Range($cellRef)
    $escaped = $text.Replace("""", """""")
    $rng.Formula = "=""" + $escaped + """"
    $rng.Copy() | Out-Null
    $rng.PasteSpecial(-4163)
}

$ws.Range("D2").Value = "30.247.82"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "1.992.62"
$ws.Range("E3").Value = "  +6.10%  "
Set-TextValue "D4" "1.002"
$ws.Range("E4").Value = "  -0.17%  "
Set-TextValue "D5" "324.66"
$ws.Range("E5").Value = "  +1.42%  "
Set-TextValue "D6" "1.001"
Set-TextValue "D7" "0.5102"
$ws.Range("E7").Value = "  +1.28%  "
Set-TextValue "D8" "0.4166"
$ws.Range("E8").Value = "  +5.28%  "
Set-TextValue "D9" "0.08714"
$ws.Range("E9").Value = "  +6.13%  "
$ws.Range("E10").Value = "  +3.59%  "
Set-TextValue "D11" "42.71"
$ws.Range("E11").Value = "  +1.57%  "
Set-TextValue "D12" "24.17"
$ws.Range("E12").Value = "  +2.74%  "
$ws.Range("D13").Value = "1.992.27"
$ws.Range("E13").Value = "  +6.13%  "
Set-TextValue "D14" "6.487"
$ws.Range("E14").Value = "  +3.06%  "
Set-TextValue "D15" "7.394"
$ws.Range("E15").Value = "  +2.73%  "
Set-TextValue "D16" "1.002"
$ws.Range("E16").Value = "  -0.10%  "
Set-TextValue "D17" "94.02"
$ws.Range("E17").Value = "  +2.43%  "
Set-TextValue "D18" "0.00001115"
$ws.Range("E18").Value = "  +2.64%  "
Set-TextValue "D19" "0.06552"
$ws.Range("E19").Value = "  +1.29%  "
Set-TextValue "D20" "18.89"
$ws.Range("E20").Value = "  +4.30%  "
Set-TextValue "D21" "1.000"
$ws.Range("E21").Value = "  -0.21%  "
Set-TextValue "D22" "6.083"
$ws.Range("E22").Value = "  +4.20%  "
$ws.Range("D23").Value = "30.304.06"
$ws.Range("E23").Value = "  +0.78%  "
Set-TextValue "D24" "11.57"
$ws.Range("E24").Value = "  +3.69%  "
Set-TextValue "D25" "2.206"
$ws.Range("E25").Value = "  +2.32%  "
$ws.Range("D26").Value = "2.221.08"
$ws.Range("E26").Value = "  +6.07%  "
$ws.Range("E27").Value = "  +6.82%  "
Set-TextValue "D28" "163.27"
$ws.Range("E28").Value = "  +1.15%  "
Set-TextValue "D29" "2.380"
$ws.Range("E29").Value = "  +6.13%  "
Set-TextValue "D30" "130.82"
$ws.Range("E30").Value = "  +2.59%  "
Set-TextValue "D31" "1.134"
$ws.Range("E31").Value = "  +5.35%  "
Set-TextValue "D32" "0.1054"
$ws.Range("E32").Value = "  +1.81%  "
Set-TextValue "D33" "6.066"
$ws.Range("E33").Value = "  +2.27%  "
$ws.Range("E34").Value = "  +3.22%  "
Set-TextValue "D35" "1.316"
$ws.Range("E35").Value = "  +12.38%  "
Set-TextValue "D36" "0.02480"
$ws.Range("E36").Value = "  +2.11%  "
Set-TextValue "D37" "5.386"
$ws.Range("E37").Value = "  +2.01%  "
Set-TextValue "D38" "0.06513"
$ws.Range("E38").Value = "  +2.45%  "
Set-TextValue "D39" "0.2188"
$ws.Range("E39").Value = "  +2.48%  "
Set-TextValue "D40" "8.920"
$ws.Range("E40").Value = "  +4.88%  "
Set-TextValue "D41" "0.6571"
$ws.Range("E41").Value = "  +4.36%  "
$ws.Range("E42").Value = "  +4.39%  "
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("E44").Value = "  +2.92%  "
Set-TextValue "D45" "0.6129"
Set-TextValue "D46" "2.192"
$ws.Range("E46").Value = "  +4.59%  "
Set-TextValue "D47" "3.660"
$ws.Range("E47").Value = "  +0.73%  "
Set-TextValue "D48" "124.32"
$ws.Range("E48").Value = "  +1.70%  "
$ws.Range("E49").Value = "  +1.39%  "
Set-TextValue "D50" "79.58"
$ws.Range("E50").Value = "  +2.77%  "
Set-TextValue "D51" "0.06878"
$ws.Range("E51").Value = "  +2.15%  "
